$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCGLT")

# Calibrate the base Max Growth Amount (MW) used by most resource rows from 800 -> 3500
$ws.Range("C2").Value = 3500

# Hydro's Max Growth Amount (row 18) is de-linked from the base value and pinned to 0
$ws.Range("C18").Formula = 0

# Restore the active cell selection on the MCGLT sheet
$ws.Activate()
$ws.Range("C3").Select()
